$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep references to the existing bodies that just move down to a different row
# (their text does not change, only their position in the sheet). These are read
# from their ORIGINAL (pre-insert) addresses, before row 13 is inserted below.
$metodoBody = $ws.Range("B19").Value2        # currently under "Critério:" (A19)
$criterioBody = $ws.Range("B20").Value2      # currently under "Norma de recuperação:" (A20)
$naoHaRecuperacao = $ws.Range("B21").Value2  # currently under "Bibliografia:" (A21)

# 1. Insert a new blank row at position 13 - this shifts rows 13-23 down to 14-24
#    (formatting/row-heights move along automatically), matching the new dimension A1:C24.
$ws.Rows("13:13").Insert()

# 2. Objetivos: body (row 10) now holds the real objectives text instead of the professor name.
$ws.Range("B10").Value = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("C10").Value = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."

# 3. The new blank row 13 (under "Docentes responsaveis:") gets the professor's name.
#    Clear the leftover A13 formatting from the row insert (A13 should stay empty/unset)
#    and copy the B/C number formats from row 14 so B13/C13 end up with the normal
#    body-text styles instead of the row-above's style.
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# 4. Programa resumido: (row 14) gets its real short-syllabus text.
$ws.Range("B14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("C14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."

# 5. Programa: (row 16) gets its real syllabus text.
$ws.Range("B16").Value = "Assuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("C16").Value = "Assuntos Técnicos específicos relacionados com o tema do projeto."

# 6. Método: (row 19, was row 18 "Método:" + professor name by mistake) gets the
#    real method description that used to sit one row below (under "Critério:").
$ws.Range("B19").Value = $metodoBody
$ws.Range("C19").Value = $metodoBody

# 7. Critério: (row 20) gets the grading-criteria text that used to sit one row below
#    (under "Norma de recuperação:").
$ws.Range("B20").Value = $criterioBody
$ws.Range("C20").Value = $criterioBody

# 8. Norma de recuperação: (row 21) gets "Não há recuperação", that used to sit one
#    row below (under "Bibliografia:").
$ws.Range("B21").Value = $naoHaRecuperacao
$ws.Range("C21").Value = $naoHaRecuperacao

# 9. Bibliografia: (row 22) gets the real bibliography text.
$ws.Range("B22").Value = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
$ws.Range("C22").Value = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."

Write-Output "done"
